$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "ZM" row (Gompertz calculation for uZM) below the existing
# beta_1 / beta_2 rows for SU, MU, R.
$ws.Range("A5").Value = "ZM"
$ws.Range("B5").Value = 0.0027550000000000001
$ws.Range("C5").Value = 0.031449999999999999

# Match the author's final selection on the sheet.
$ws.Range("C10").Select()
